# Automatic update of files.
# Rotates the A / I / Q / R / AC values among rows 2, 3 and 4:
#   new row2 <- old row3
#   new row3 <- old row4
#   new row4 <- old row2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the current ("old") values before we overwrite anything ---
$oldA2 = $ws.Range("A2").Value2
$oldI2 = $ws.Range("I2").Value2
$oldQ2 = $ws.Range("Q2").Value2
$oldR2 = $ws.Range("R2").Value2
$oldAC2 = $ws.Range("AC2").Value2

$oldA3 = $ws.Range("A3").Value2
$oldI3 = $ws.Range("I3").Value2
$oldQ3 = $ws.Range("Q3").Value2
$oldR3 = $ws.Range("R3").Value2
$oldAC3 = $ws.Range("AC3").Value2

$oldA4 = $ws.Range("A4").Value2
$oldI4 = $ws.Range("I4").Value2
$oldQ4 = $ws.Range("Q4").Value2
$oldR4 = $ws.Range("R4").Value2
$oldAC4 = $ws.Range("AC4").Value2

# The "Antal" column (I) is stored as text in this workbook, so force the
# text number format before writing numeric-looking strings into it.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"

# --- row 2 <- old row 3 ---
$ws.Range("A2").Value = $oldA3
$ws.Range("I2").Value = $oldI3
$ws.Range("Q2").Value = $oldQ3
$ws.Range("R2").Value = $oldR3
if ($oldAC3 -eq $null) {
    $ws.Range("AC2").Value = ""
} else {
    $ws.Range("AC2").Value = $oldAC3
}

# --- row 3 <- old row 4 ---
$ws.Range("A3").Value = $oldA4
$ws.Range("I3").Value = $oldI4
$ws.Range("Q3").Value = $oldQ4
$ws.Range("R3").Value = $oldR4
if ($oldAC4 -eq $null) {
    $ws.Range("AC3").Value = ""
} else {
    $ws.Range("AC3").Value = $oldAC4
}

# --- row 4 <- old row 2 ---
$ws.Range("A4").Value = $oldA2
$ws.Range("I4").Value = $oldI2
$ws.Range("Q4").Value = $oldQ2
$ws.Range("R4").Value = $oldR2
if ($oldAC2 -eq $null) {
    $ws.Range("AC4").Value = ""
} else {
    $ws.Range("AC4").Value = $oldAC2
}
